$d = $word.ActiveDocument

# 1) Remove the _GoBack bookmark currently sitting after "Trừ:" in the
#    "Hình thức chấm điểm" bullet list.
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# 2) Fix the "Tổng kết giữa năm ... nhất trị giá 300k." sentence so it
#    reads "... nhất." (drop " trị giá 300k") and ends up split across
#    three runs: "...nhấ" / "t" / "." — matching the target OOXML.
$rng = $d.Content
$found = $rng.Find.Execute("nhất trị giá 300k.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found) {
    $matchStart = $rng.Start
    $matchEnd = $rng.End

    # The matched text is "nhất trị giá 300k." -- remove everything
    # between "nhấ" + "t" (keep both) and the final "." (keep it too),
    # i.e. delete " trị giá 300k" which sits right after the final "t"
    # of "nhất" and right before the trailing ".".
    $tStart = $matchStart + 3   # start of the final "t" in "nhất"
    $tEnd = $tStart + 1         # end of that "t"
    $periodStart = $matchEnd - 1  # start of the trailing "."

    $middle = $d.Range($tEnd, $periodStart)
    $middle.Delete()

    # Force the single merged run to split into three runs with no
    # leftover formatting markup, by bracketing the lone "t" character
    # with a temporary bookmark (bookmarks force run boundaries but
    # leave no <w:rPr/> residue) and then removing the bookmark again.
    $tRange = $d.Range($tStart, $tEnd)
    $d.Bookmarks.Add("__tmp_split__", $tRange)
    $d.Bookmarks("__tmp_split__").Delete()
}

# 3) Re-add the _GoBack bookmark to the empty paragraph at the very end
#    of the document (its new home after the edit).
$lastPara = $d.Paragraphs.Last
$endRange = $lastPara.Range
$endRange.Collapse(0)
$d.Bookmarks.Add("_GoBack", $endRange)
